$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 6.345999999999999
$ws.Range("D12").Value = -7.269
$ws.Range("E15").Value = 16.158
$ws.Range("B18").Value = 5.166
$ws.Range("B20").Value = 6.992
$ws.Range("D26").Value = -8.062000000000001
$ws.Range("B27").Value = 6.105
$ws.Range("D27").Value = -8.251999999999999
$ws.Range("D29").Value = -7.297999999999999
$ws.Range("D37").Value = -7.822
$ws.Range("D38").Value = -7.355000000000001
$ws.Range("E38").Value = 16.658
$ws.Range("E44").Value = 16.874
$ws.Range("D51").Value = -8.401
$ws.Range("E51").Value = 16.634
$ws.Range("D55").Value = -8.129000000000001
$ws.Range("E57").Value = 16.532
$ws.Range("E63").Value = 17.601
$ws.Range("B69").Value = 5.992
$ws.Range("D69").Value = -7.438000000000001
$ws.Range("D70").Value = -7.175999999999999
$ws.Range("E70").Value = 17.792
$ws.Range("B76").Value = 6.308
$ws.Range("B82").Value = 5.366
$ws.Range("D83").Value = -8.373000000000001
$ws.Range("E99").Value = 16.643
$ws.Range("D102").Value = -7.833
